$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (insert new Jan_2026 column, shift Dec_2025/Nov_2025, drop Oct_2025)
$ws.Cells.Item(1,1).Value() = "ISIN"
$ws.Cells.Item(1,2).Value() = "Stock Name"
$ws.Cells.Item(1,3).Value() = "Mutual Fund"
$ws.Cells.Item(1,4).Value() = "Jan_2026"
$ws.Cells.Item(1,5).Value() = "Dec_2025"
$ws.Cells.Item(1,6).Value() = "Nov_2025"
$ws.Cells.Item(1,7).Value() = "MoM"
$ws.Cells.Item(1,8).Value() = "QoQ"

# Update data rows 2-20 with refreshed holdings data
$ws.Cells.Item(2,1).Value() = "INE781S01027"
$ws.Cells.Item(2,2).Value() = "Ventive Hospitality Limited"
$ws.Cells.Item(2,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(2,4).Value() = 9.802079
$ws.Cells.Item(2,5).Value() = 10.033282
$ws.Cells.Item(2,6).Value() = 9.688181
$ws.Cells.Item(2,7).Value() = -0.2312029999999989
$ws.Cells.Item(2,8).Value() = 0.1138980000000007
$ws.Cells.Item(3,1).Value() = "INE180C01042"
$ws.Cells.Item(3,2).Value() = "Capri Global Capital Limited"
$ws.Cells.Item(3,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(3,4).Value() = 9.526872
$ws.Cells.Item(3,5).Value() = 9.512082
$ws.Cells.Item(3,6).Value() = 9.425583
$ws.Cells.Item(3,7).Value() = 0.01478999999999964
$ws.Cells.Item(3,8).Value() = 0.1012889999999995
$ws.Cells.Item(4,1).Value() = "INE768C01028"
$ws.Cells.Item(4,2).Value() = "Zydus Wellness Ltd"
$ws.Cells.Item(4,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(4,4).Value() = 8.403816
$ws.Cells.Item(4,5).Value() = 9.965963
$ws.Cells.Item(4,6).Value() = 9.087308
$ws.Cells.Item(4,7).Value() = -1.562147
$ws.Cells.Item(4,8).Value() = -0.6834919999999993
$ws.Cells.Item(5,1).Value() = "INE016A01026"
$ws.Cells.Item(5,2).Value() = "Dabur India Limited"
$ws.Cells.Item(5,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(5,4).Value() = 7.315448
$ws.Cells.Item(5,5).Value() = 3.949244
$ws.Cells.Item(5,6).Value() = 3.907306
$ws.Cells.Item(5,7).Value() = 3.366204
$ws.Cells.Item(5,8).Value() = 3.408142
$ws.Cells.Item(6,1).Value() = "INE917I01010"
$ws.Cells.Item(6,2).Value() = "Bajaj Auto Limited"
$ws.Cells.Item(6,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(6,4).Value() = 7.003441
$ws.Cells.Item(6,5).Value() = 6.379967
$ws.Cells.Item(6,6).Value() = 5.966638
$ws.Cells.Item(6,7).Value() = 0.6234739999999999
$ws.Cells.Item(6,8).Value() = 1.036803
$ws.Cells.Item(7,1).Value() = "INE804L01022"
$ws.Cells.Item(7,2).Value() = "Medplus Health Services Limited"
$ws.Cells.Item(7,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(7,4).Value() = 5.84704
$ws.Cells.Item(7,5).Value() = 5.528791
$ws.Cells.Item(7,6).Value() = 5.334075
$ws.Cells.Item(7,7).Value() = 0.3182489999999998
$ws.Cells.Item(7,8).Value() = 0.5129649999999994
$ws.Cells.Item(8,1).Value() = "INE406A01037"
$ws.Cells.Item(8,2).Value() = "Aurobindo Pharma Limited"
$ws.Cells.Item(8,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(8,4).Value() = 5.578617
$ws.Cells.Item(8,5).Value() = 3.282408
$ws.Cells.Item(8,6).Value() = 3.277698
$ws.Cells.Item(8,7).Value() = 2.296209
$ws.Cells.Item(8,8).Value() = 2.300919
$ws.Cells.Item(9,1).Value() = "INE04TZ01018"
$ws.Cells.Item(9,2).Value() = "ETHOS LIMITED"
$ws.Cells.Item(9,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(9,4).Value() = 5.232357
$ws.Cells.Item(9,5).Value() = 5.71073
$ws.Cells.Item(9,6).Value() = 2.769367
$ws.Cells.Item(9,7).Value() = -0.4783729999999995
$ws.Cells.Item(9,8).Value() = 2.46299
$ws.Cells.Item(10,1).Value() = "INE192A01025"
$ws.Cells.Item(10,2).Value() = "Tata Consumer Products Ltd"
$ws.Cells.Item(10,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(10,4).Value() = 3.671604
$ws.Cells.Item(10,5).Value() = 0.493396
$ws.Cells.Item(10,6).Value() = 0
$ws.Cells.Item(10,7).Value() = 3.178208
$ws.Cells.Item(10,8).Value() = 3.671604
$ws.Cells.Item(11,1).Value() = "INE179A01014"
$ws.Cells.Item(11,2).Value() = "Procter & Gamble Hygiene & Health Care Limited"
$ws.Cells.Item(11,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(11,4).Value() = 3.426913
$ws.Cells.Item(11,5).Value() = 3.535301
$ws.Cells.Item(11,6).Value() = 3.348036
$ws.Cells.Item(11,7).Value() = -0.1083880000000002
$ws.Cells.Item(11,8).Value() = 0.07887699999999986
$ws.Cells.Item(12,1).Value() = "INE01A001028"
$ws.Cells.Item(12,2).Value() = "Stanley Lifestyles Limited"
$ws.Cells.Item(12,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(12,4).Value() = 1.90172
$ws.Cells.Item(12,5).Value() = 3.168372
$ws.Cells.Item(12,6).Value() = 4.822661
$ws.Cells.Item(12,7).Value() = -1.266652
$ws.Cells.Item(12,8).Value() = -2.920941
$ws.Cells.Item(13,1).Value() = "INE018E01016"
$ws.Cells.Item(13,2).Value() = "SBI Cards & Payment Services Ltd"
$ws.Cells.Item(13,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(13,4).Value() = 0.990131
$ws.Cells.Item(13,5).Value() = 1.059535
$ws.Cells.Item(13,6).Value() = 1.04217
$ws.Cells.Item(13,7).Value() = -0.06940399999999991
$ws.Cells.Item(13,8).Value() = -0.05203900000000006
$ws.Cells.Item(14,1).Value() = "INE090A01021"
$ws.Cells.Item(14,2).Value() = "ICICI Bank Limited"
$ws.Cells.Item(14,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(14,4).Value() = 0.418628
$ws.Cells.Item(14,5).Value() = 0
$ws.Cells.Item(14,6).Value() = 0
$ws.Cells.Item(14,7).Value() = 0.418628
$ws.Cells.Item(14,8).Value() = 0.418628
$ws.Cells.Item(15,1).Value() = "INE196A01026"
$ws.Cells.Item(15,2).Value() = "Marico Limited"
$ws.Cells.Item(15,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(15,4).Value() = 0
$ws.Cells.Item(15,5).Value() = 6.171663
$ws.Cells.Item(15,6).Value() = 0
$ws.Cells.Item(15,7).Value() = -6.171663
$ws.Cells.Item(15,8).Value() = 0
$ws.Cells.Item(16,1).Value() = "INE484J01027"
$ws.Cells.Item(16,2).Value() = "Godrej Properties Limited"
$ws.Cells.Item(16,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(16,4).Value() = 0
$ws.Cells.Item(16,5).Value() = 6.542531
$ws.Cells.Item(16,6).Value() = 3.848956
$ws.Cells.Item(16,7).Value() = -6.542531
$ws.Cells.Item(16,8).Value() = -3.848956
$ws.Cells.Item(17,1).Value() = "INE669C01036"
$ws.Cells.Item(17,2).Value() = "Tech Mahindra Limited"
$ws.Cells.Item(17,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(17,4).Value() = 0
$ws.Cells.Item(17,5).Value() = 0.272518
$ws.Cells.Item(17,6).Value() = 0
$ws.Cells.Item(17,7).Value() = -0.272518
$ws.Cells.Item(17,8).Value() = 0
$ws.Cells.Item(18,1).Value() = "INE686F01025"
$ws.Cells.Item(18,2).Value() = "UNITED BREWERIES LIMITED"
$ws.Cells.Item(18,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(18,4).Value() = 0
$ws.Cells.Item(18,5).Value() = 0
$ws.Cells.Item(18,6).Value() = 2.664613
$ws.Cells.Item(18,7).Value() = 0
$ws.Cells.Item(18,8).Value() = -2.664613
$ws.Cells.Item(19,1).Value() = "INE854D01024"
$ws.Cells.Item(19,2).Value() = "United Spirits Limited"
$ws.Cells.Item(19,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(19,4).Value() = 0
$ws.Cells.Item(19,5).Value() = 5.232451
$ws.Cells.Item(19,6).Value() = 5.066381
$ws.Cells.Item(19,7).Value() = -5.232451
$ws.Cells.Item(19,8).Value() = -5.066381
$ws.Cells.Item(20,1).Value() = "INE202B01038"
$ws.Cells.Item(20,2).Value() = "Piramal Finance Ltd"
$ws.Cells.Item(20,3).Value() = "quant Consumption Fund"
$ws.Cells.Item(20,4).Value() = 0
$ws.Cells.Item(20,5).Value() = 0
$ws.Cells.Item(20,6).Value() = 3.103411
$ws.Cells.Item(20,7).Value() = 0
$ws.Cells.Item(20,8).Value() = -3.103411

# Remove the now-obsolete trailing rows (21-23) that fell out of the comparison window
$ws.Range("A21:H23").ClearContents()
